$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Uncut_Sheet_1"
